# Refresh the cryptocurrency price/volume snapshot in the "cryptos" workbook,
# mirroring the data pulled by the scheduled GitHub Actions job.
# For every affected row, Price (column D) and/or Volume(1h) (column E) are
# updated to their newly scraped values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "70.916.12"; E = "  +2.19%  " }
    @{ Row = 3; D = "3.820.35"; E = "  +0.90%  " }
    @{ Row = 4; D = $null; E = "  +0.03%  " }
    @{ Row = 5; D = "669.71"; E = "  +7.21%  " }
    @{ Row = 6; D = "169.74"; E = "  +2.38%  " }
    @{ Row = 7; D = "3.818.52"; E = "  +0.92%  " }
    @{ Row = 8; D = $null; E = "  -0.01%  " }
    @{ Row = 9; D = $null; E = "  +1.26%  " }
    @{ Row = 10; D = $null; E = "  +0.70%  " }
    @{ Row = 11; D = $null; E = "  +2.57%  " }
    @{ Row = 12; D = "6.98"; E = "  +4.27%  " }
    @{ Row = 13; D = $null; E = "  -0.40%  " }
    @{ Row = 14; D = "36.12"; E = "  +1.18%  " }
    @{ Row = 15; D = "4.467.11"; E = "  +1.06%  " }
    @{ Row = 16; D = "3.823.04"; E = "  +1.32%  " }
    @{ Row = 17; D = "70.859.23"; E = "  +2.19%  " }
    @{ Row = 18; D = "17.79"; E = "  +0.51%  " }
    @{ Row = 19; D = $null; E = "  +21.98%  " }
    @{ Row = 20; D = $null; E = "  +1.04%  " }
    @{ Row = 21; D = $null; E = "  +0.67%  " }
    @{ Row = 22; D = "476.32"; E = "  +1.66%  " }
    @{ Row = 23; D = $null; E = "  +1.87%  " }
    @{ Row = 24; D = "83.33"; E = "  +0.02%  " }
    @{ Row = 25; D = "0.0000146"; E = "  -2.33%  " }
    @{ Row = 26; D = "12.26"; E = "  +1.83%  " }
    @{ Row = 27; D = $null; E = "  +3.56%  " }
    @{ Row = 28; D = "2.13"; E = "  -1.44%  " }
    @{ Row = 29; D = $null; E = "  -0.01%  " }
    @{ Row = 30; D = "3.973.96"; E = "  +0.98%  " }
    @{ Row = 31; D = "2.88"; E = "  +8.15%  " }
    @{ Row = 32; D = "2.31"; E = "  +2.94%  " }
    @{ Row = 33; D = $null; E = "  +2.37%  " }
    @{ Row = 34; D = "29.85"; E = "  +3.53%  " }
    @{ Row = 35; D = $null; E = "  +6.30%  " }
    @{ Row = 36; D = "9.22"; E = "  +2.16%  " }
    @{ Row = 37; D = "3.777.92"; E = "  +1.11%  " }
    @{ Row = 38; D = $null; E = "  -0.12%  " }
    @{ Row = 39; D = $null; E = "  +0.44%  " }
    @{ Row = 40; D = "3.46"; E = "  +1.59%  " }
    @{ Row = 41; D = "6.01"; E = "  +3.27%  " }
    @{ Row = 42; D = "0.968"; E = "  -0.17%  " }
    @{ Row = 43; D = "1.00"; E = "  +0.10%  " }
    @{ Row = 44; D = "2.12"; E = "  +9.95%  " }
    @{ Row = 45; D = $null; E = "  -0.02%  " }
    @{ Row = 46; D = $null; E = "  +5.06%  " }
    @{ Row = 47; D = "157.74"; E = "  +3.61%  " }
    @{ Row = 48; D = "48.10"; E = "  +2.86%  " }
    @{ Row = 49; D = $null; E = "  +0.84%  " }
    @{ Row = 50; D = "1.42"; E = "  +3.88%  " }
    @{ Row = 51; D = "8.53"; E = "  +1.17%  " }
)

foreach ($update in $updates) {
    $row = $update.Row

    if ($null -ne $update.D) {
        $cell = $ws.Range("D$row")
        # Force a text number format first so Excel keeps these values
        # (prices such as "70.916.12" or "1.00") as literal text instead of
        # silently re-interpreting them as numbers/dates.
        $cell.NumberFormat = "@"
        $cell.Value = $update.D
    }

    if ($null -ne $update.E) {
        $cell = $ws.Range("E$row")
        $cell.NumberFormat = "@"
        $cell.Value = $update.E
    }
}
